$wb = $excel.ActiveWorkbook

# ALC!row33
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 149.4
$ws.Cells.Item(33, 9).Value = 149.4
$ws.Cells.Item(33, 11).Value = 149.4
$ws.Cells.Item(33, 13).Value = 79.59999999999999

# ALC!row74
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 4750
$ws.Cells.Item(74, 10).Value = 5000
$ws.Cells.Item(74, 12).Value = 5000
$ws.Cells.Item(74, 14).Value = -6872

# ALC!row77
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(77, 8).Value = 4750
$ws.Cells.Item(77, 10).Value = 5000
$ws.Cells.Item(77, 12).Value = 25000
$ws.Cells.Item(77, 14).Value = -34360

# ARM!row8
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(8, 8).Value = 5.714286
$ws.Cells.Item(8, 9).Value = 5.714286
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 5.714286
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = 138.285714
$ws.Cells.Item(8, 14).Value = $null

# ARM!row11
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(11, 8).Value = 1005
$ws.Cells.Item(11, 10).Value = 2000
$ws.Cells.Item(11, 12).Value = 2000
$ws.Cells.Item(11, 14).Value = -2288

# ARM!row12
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(12, 8).Value = 7000
$ws.Cells.Item(12, 10).Value = 7000
$ws.Cells.Item(12, 12).Value = 7000
$ws.Cells.Item(12, 14).Value = -7346

# ARM!row13
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(13, 8).Value = 10
$ws.Cells.Item(13, 9).Value = 10
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 10
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = 134
$ws.Cells.Item(13, 14).Value = $null

# ARM!row14
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(14, 8).Value = 4333.3335
$ws.Cells.Item(14, 9).Value = 6000
$ws.Cells.Item(14, 10).Value = 1000
$ws.Cells.Item(14, 11).Value = 6000
$ws.Cells.Item(14, 12).Value = 1000
$ws.Cells.Item(14, 13).Value = -5825
$ws.Cells.Item(14, 14).Value = -1350

# ARM!row15
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(15, 8).Value = 3125
$ws.Cells.Item(15, 10).Value = 3125
$ws.Cells.Item(15, 12).Value = 3125
$ws.Cells.Item(15, 14).Value = -3825

# ARM!row17
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(17, 8).Value = 53.5
$ws.Cells.Item(17, 9).Value = 53.5
$ws.Cells.Item(17, 11).Value = 53.5
$ws.Cells.Item(17, 13).Value = 119.5

# ARM!row18
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 14).Value = $null

# ARM!row19
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(19, 8).Value = 10
$ws.Cells.Item(19, 10).Value = 10
$ws.Cells.Item(19, 12).Value = 10
$ws.Cells.Item(19, 14).Value = -468

# ARM!row22
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 13).Value = $null

# ARM!row25
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(25, 8).Value = 4245.3335
$ws.Cells.Item(25, 9).Value = 3094.4
$ws.Cells.Item(25, 10).Value = 10000
$ws.Cells.Item(25, 11).Value = 3094.4
$ws.Cells.Item(25, 12).Value = 10000
$ws.Cells.Item(25, 13).Value = -2692.4
$ws.Cells.Item(25, 14).Value = -10804

# ARM!row33
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(33, 8).Value = 26
$ws.Cells.Item(33, 9).Value = 26
$ws.Cells.Item(33, 11).Value = 26
$ws.Cells.Item(33, 13).Value = 303

# ARM!row36
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(36, 8).Value = 5492.7144
$ws.Cells.Item(36, 9).Value = 5289.8
$ws.Cells.Item(36, 11).Value = 5289.8
$ws.Cells.Item(36, 13).Value = -4943.8

# ARM!row63
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 5775
$ws.Cells.Item(63, 10).Value = 5000
$ws.Cells.Item(63, 12).Value = 5000
$ws.Cells.Item(63, 14).Value = -6372

# ARM!row66
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 5775
$ws.Cells.Item(66, 10).Value = 5000
$ws.Cells.Item(66, 12).Value = 25000
$ws.Cells.Item(66, 14).Value = -31864

# ARM!row97
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 1110
$ws.Cells.Item(97, 9).Value = 1110
$ws.Cells.Item(97, 11).Value = 1110
$ws.Cells.Item(97, 13).Value = -614

# ARM!row102
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 773.8
$ws.Cells.Item(102, 9).Value = 642.5
$ws.Cells.Item(102, 11).Value = 642.5
$ws.Cells.Item(102, 13).Value = 979.5

# ARM!row110
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 586.7143
$ws.Cells.Item(110, 9).Value = 586.7143
$ws.Cells.Item(110, 11).Value = 586.7143
$ws.Cells.Item(110, 13).Value = 1458.2857

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 5000
$ws.Cells.Item(122, 9).Value = 5000
$ws.Cells.Item(122, 11).Value = 15000
$ws.Cells.Item(122, 13).Value = -12550

# BSM!row36
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(36, 8).Value = 3805.4
$ws.Cells.Item(36, 9).Value = 3805.4
$ws.Cells.Item(36, 11).Value = 3805.4
$ws.Cells.Item(36, 13).Value = -3271.4

# BSM!row94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1906.7142
$ws.Cells.Item(94, 9).Value = 1906.7142
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 1906.7142
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 13).Value = -1455.7142
$ws.Cells.Item(94, 14).Value = $null

# BSM!row99
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1302.4286
$ws.Cells.Item(99, 9).Value = 1311.6666
$ws.Cells.Item(99, 11).Value = 1311.6666
$ws.Cells.Item(99, 13).Value = 186.3334

# BSM!row107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1064.091
$ws.Cells.Item(107, 9).Value = 1020.5
$ws.Cells.Item(107, 10).Value = 1500
$ws.Cells.Item(107, 11).Value = 1020.5
$ws.Cells.Item(107, 12).Value = 1500
$ws.Cells.Item(107, 13).Value = 899.5
$ws.Cells.Item(107, 14).Value = -5340

# BSM!row108
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 14).Value = $null

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6593.933
$ws.Cells.Item(31, 9).Value = 1354.8572
$ws.Cells.Item(31, 11).Value = 1354.8572
$ws.Cells.Item(31, 13).Value = -1059.8572

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 6593.933
$ws.Cells.Item(34, 9).Value = 1354.8572
$ws.Cells.Item(34, 11).Value = 1354.8572
$ws.Cells.Item(34, 13).Value = -1152.8572

# CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 17998.5
$ws.Cells.Item(58, 10).Value = 17998.5
$ws.Cells.Item(58, 12).Value = 17998.5
$ws.Cells.Item(58, 14).Value = -18404.5

# CRP!row60
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(60, 8).Value = 10000
$ws.Cells.Item(60, 9).Value = 10000
$ws.Cells.Item(60, 10).Value = 0
$ws.Cells.Item(60, 11).Value = 10000
$ws.Cells.Item(60, 12).Value = 0
$ws.Cells.Item(60, 13).Value = -9489
$ws.Cells.Item(60, 14).Value = $null

# CRP!row107
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 1333.3334
$ws.Cells.Item(107, 9).Value = 2000
$ws.Cells.Item(107, 11).Value = 2000
$ws.Cells.Item(107, 13).Value = -80

# CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 17998.5
$ws.Cells.Item(136, 10).Value = 17998.5
$ws.Cells.Item(136, 12).Value = 53995.5
$ws.Cells.Item(136, 14).Value = -59095.5

# CUL!row63
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 13).Value = $null

# CUL!row66
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(66, 8).Value = 0
$ws.Cells.Item(66, 9).Value = 0
$ws.Cells.Item(66, 11).Value = 0
$ws.Cells.Item(66, 13).Value = $null

# CUL!row97
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(97, 8).Value = 76.5
$ws.Cells.Item(97, 10).Value = 76
$ws.Cells.Item(97, 12).Value = 228
$ws.Cells.Item(97, 14).Value = -1220

# GSM!row9
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(9, 8).Value = 1999
$ws.Cells.Item(9, 10).Value = 1999
$ws.Cells.Item(9, 12).Value = 1999
$ws.Cells.Item(9, 14).Value = -2339

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 549.6667
$ws.Cells.Item(7, 9).Value = 599.5
$ws.Cells.Item(7, 10).Value = 450
$ws.Cells.Item(7, 11).Value = 599.5
$ws.Cells.Item(7, 12).Value = 450
$ws.Cells.Item(7, 13).Value = -487.5
$ws.Cells.Item(7, 14).Value = -674

# LTW!row22
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 956.125
$ws.Cells.Item(22, 9).Value = 850
$ws.Cells.Item(22, 10).Value = 1133
$ws.Cells.Item(22, 11).Value = 850
$ws.Cells.Item(22, 12).Value = 1133
$ws.Cells.Item(22, 13).Value = -555
$ws.Cells.Item(22, 14).Value = -1723

# LTW!row27
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 956.125
$ws.Cells.Item(27, 9).Value = 850
$ws.Cells.Item(27, 10).Value = 1133
$ws.Cells.Item(27, 11).Value = 850
$ws.Cells.Item(27, 12).Value = 1133
$ws.Cells.Item(27, 13).Value = -743
$ws.Cells.Item(27, 14).Value = -1347

# LTW!row68
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 4045
$ws.Cells.Item(68, 9).Value = 3192.5
$ws.Cells.Item(68, 10).Value = 5750
$ws.Cells.Item(68, 11).Value = 3192.5
$ws.Cells.Item(68, 12).Value = 5750
$ws.Cells.Item(68, 13).Value = -2443.5
$ws.Cells.Item(68, 14).Value = -7248

# LTW!row71
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 4045
$ws.Cells.Item(71, 9).Value = 3192.5
$ws.Cells.Item(71, 10).Value = 5750
$ws.Cells.Item(71, 11).Value = 15962.5
$ws.Cells.Item(71, 12).Value = 28750
$ws.Cells.Item(71, 13).Value = -12218.5
$ws.Cells.Item(71, 14).Value = -36238

# LTW!row98
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 14).Value = $null

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 549.6667
$ws.Cells.Item(126, 9).Value = 599.5
$ws.Cells.Item(126, 10).Value = 450
$ws.Cells.Item(126, 11).Value = 1798.5
$ws.Cells.Item(126, 12).Value = 1350
$ws.Cells.Item(126, 13).Value = 671.5
$ws.Cells.Item(126, 14).Value = -6290

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 10234.75
$ws.Cells.Item(132, 9).Value = 10234.75
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 30704.25
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -28174.25
$ws.Cells.Item(132, 14).Value = $null

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 15999
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 15999
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).Value = 47997
$ws.Cells.Item(136, 13).Value = $null
$ws.Cells.Item(136, 14).Value = -53097
